# Refresh the "want to go" counts (column F) and, on two rows, the minimum
# ticket price (column G) to match the live data pulled at commit 456a3b4.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 653
$ws.Range("F7").Value = 2167
$ws.Range("F8").Value = 904
$ws.Range("F9").Value = 859
$ws.Range("F12").Value = 429
$ws.Range("F13").Value = 325
$ws.Range("F14").Value = 103
$ws.Range("F15").Value = 1101
$ws.Range("F18").Value = 1806
$ws.Range("F26").Value = 530
$ws.Range("F27").Value = 357
$ws.Range("F29").Value = 426
$ws.Range("F30").Value = 2504
$ws.Range("G30").Value = 59.9
$ws.Range("F31").Value = 387
$ws.Range("F34").Value = 607
$ws.Range("F35").Value = 483
$ws.Range("F37").Value = 933
$ws.Range("F40").Value = 519
$ws.Range("F41").Value = 510

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 23
$ws.Range("F15").Value = 5
$ws.Range("F22").Value = 125
$ws.Range("F23").Value = 114

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 229
$ws.Range("F5").Value = 247
$ws.Range("F6").Value = 332

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 229
$ws.Range("F7").Value = 247
$ws.Range("F8").Value = 653
$ws.Range("F11").Value = 904
$ws.Range("F12").Value = 859
$ws.Range("F15").Value = 429
$ws.Range("F16").Value = 325
$ws.Range("F17").Value = 103
$ws.Range("F19").Value = 1101
$ws.Range("F23").Value = 332
$ws.Range("F24").Value = 1806
$ws.Range("F33").Value = 5
$ws.Range("F35").Value = 530
$ws.Range("F37").Value = 426
$ws.Range("F38").Value = 2504
$ws.Range("G38").Value = 59.9
$ws.Range("F40").Value = 607
$ws.Range("F41").Value = 483
$ws.Range("F43").Value = 933
$ws.Range("F44").Value = 114
$ws.Range("F48").Value = 519
$ws.Range("F49").Value = 510
